$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MISC"
$ws.Range("A3").Value = "RN_LIBRARIES"
